$p = $ppt.ActivePresentation

$oldDate = "2021/11/30"
$newDate = "2023/5/6"

# 1. Update the date placeholder text baked into the slide master.
for ($i = 1; $i -le $p.SlideMaster.Shapes.Count; $i++) {
    $shp = $p.SlideMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# 2. Update the date placeholder text baked into every slide layout.
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 3. Fix the wording on slide 17: "a virtual destructor" -> "a destructor".
$s = $p.Slides.Item(17)
$shp = $s.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "If a destructor is not virtual, only the destructor of the base class is executed in the follow examples."
